$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-13 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-14 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("97-32=65", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=59", 2) | Out-Null
$d.Content.Find.Execute("25+8=33", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=5", 2) | Out-Null
$d.Content.Find.Execute("32+4=36", $true, $false, $false, $false, $false, $true, 1, $false, "95-0=95", 2) | Out-Null
$d.Content.Find.Execute("92-60=32", $true, $false, $false, $false, $false, $true, 1, $false, "20-15=5", 2) | Out-Null
$d.Content.Find.Execute("21+60=81", $true, $false, $false, $false, $false, $true, 1, $false, "77+6=83", 2) | Out-Null
$d.Content.Find.Execute("4+69=73", $true, $false, $false, $false, $false, $true, 1, $false, "54-28=26", 2) | Out-Null
$d.Content.Find.Execute("69-37=32", $true, $false, $false, $false, $false, $true, 1, $false, "48+23=71", 2) | Out-Null
$d.Content.Find.Execute("58+20=78", $true, $false, $false, $false, $false, $true, 1, $false, "26+49=75", 2) | Out-Null
$d.Content.Find.Execute("21-16=5", $true, $false, $false, $false, $false, $true, 1, $false, "18+28=46", 2) | Out-Null
$d.Content.Find.Execute("88-19=69", $true, $false, $false, $false, $false, $true, 1, $false, "53+31=84", 2) | Out-Null
$d.Content.Find.Execute("95-26=69", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=24", 2) | Out-Null
$d.Content.Find.Execute("69-61=8", $true, $false, $false, $false, $false, $true, 1, $false, "9+30=39", 2) | Out-Null
$d.Content.Find.Execute("11+1=12", $true, $false, $false, $false, $false, $true, 1, $false, "59+4=63", 2) | Out-Null
$d.Content.Find.Execute("17+9=26", $true, $false, $false, $false, $false, $true, 1, $false, "82-34=48", 2) | Out-Null
$d.Content.Find.Execute("29-11=18", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=22", 2) | Out-Null
$d.Content.Find.Execute("15+32=47", $true, $false, $false, $false, $false, $true, 1, $false, "88-5=83", 2) | Out-Null
$d.Content.Find.Execute("45-15=30", $true, $false, $false, $false, $false, $true, 1, $false, "83-17=66", 2) | Out-Null
$d.Content.Find.Execute("13+36=49", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=72", 2) | Out-Null
$d.Content.Find.Execute("82-13=69", $true, $false, $false, $false, $false, $true, 1, $false, "48-33=15", 2) | Out-Null
$d.Content.Find.Execute("4+82=86", $true, $false, $false, $false, $false, $true, 1, $false, "50+3=53", 2) | Out-Null
$d.Content.Find.Execute("13+12=25", $true, $false, $false, $false, $false, $true, 1, $false, "6+33=39", 2) | Out-Null
$d.Content.Find.Execute("84-45=39", $true, $false, $false, $false, $false, $true, 1, $false, "13+49=62", 2) | Out-Null
$d.Content.Find.Execute("83-61=22", $true, $false, $false, $false, $false, $true, 1, $false, "37+14=51", 2) | Out-Null
$d.Content.Find.Execute("29-24=5", $true, $false, $false, $false, $false, $true, 1, $false, "6+91=97", 2) | Out-Null
$d.Content.Find.Execute("41+50=91", $true, $false, $false, $false, $false, $true, 1, $false, "29+6=35", 2) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "62+36=98", 2) | Out-Null
$d.Content.Find.Execute("7+19=26", $true, $false, $false, $false, $false, $true, 1, $false, "6+19=25", 2) | Out-Null
$d.Content.Find.Execute("83-73=10", $true, $false, $false, $false, $false, $true, 1, $false, "69-5=64", 2) | Out-Null
$d.Content.Find.Execute("48-39=9", $true, $false, $false, $false, $false, $true, 1, $false, "73+9=82", 2) | Out-Null
$d.Content.Find.Execute("88-36=52", $true, $false, $false, $false, $false, $true, 1, $false, "58-36=22", 2) | Out-Null
$d.Content.Find.Execute("44-35=9", $true, $false, $false, $false, $false, $true, 1, $false, "55-18=37", 2) | Out-Null
$d.Content.Find.Execute("56-7=49", $true, $false, $false, $false, $false, $true, 1, $false, "36-19=17", 2) | Out-Null
$d.Content.Find.Execute("66-2=64", $true, $false, $false, $false, $false, $true, 1, $false, "42+14=56", 2) | Out-Null
$d.Content.Find.Execute("83-70=13", $true, $false, $false, $false, $false, $true, 1, $false, "19-12=7", 2) | Out-Null
$d.Content.Find.Execute("97-25=72", $true, $false, $false, $false, $false, $true, 1, $false, "55+42=97", 2) | Out-Null
$d.Content.Find.Execute("13+42=55", $true, $false, $false, $false, $false, $true, 1, $false, "60+18=78", 2) | Out-Null
$d.Content.Find.Execute("24-18=6", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=38", 2) | Out-Null
$d.Content.Find.Execute("37+42=79", $true, $false, $false, $false, $false, $true, 1, $false, "36+27=63", 2) | Out-Null
$d.Content.Find.Execute("53+5=58", $true, $false, $false, $false, $false, $true, 1, $false, "87-70=17", 2) | Out-Null
$d.Content.Find.Execute("38+41=79", $true, $false, $false, $false, $false, $true, 1, $false, "45-18=27", 2) | Out-Null
$d.Content.Find.Execute("70+23=93", $true, $false, $false, $false, $false, $true, 1, $false, "13+18=31", 2) | Out-Null
$d.Content.Find.Execute("34+28=62", $true, $false, $false, $false, $false, $true, 1, $false, "81-17=64", 2) | Out-Null
$d.Content.Find.Execute("53+34=87", $true, $false, $false, $false, $false, $true, 1, $false, "82-25=57", 2) | Out-Null
$d.Content.Find.Execute("36+12=48", $true, $false, $false, $false, $false, $true, 1, $false, "9+19=28", 2) | Out-Null
$d.Content.Find.Execute("46-4=42", $true, $false, $false, $false, $false, $true, 1, $false, "12+62=74", 2) | Out-Null
$d.Content.Find.Execute("33+24=57", $true, $false, $false, $false, $false, $true, 1, $false, "63-0=63", 2) | Out-Null
$d.Content.Find.Execute("29+63=92", $true, $false, $false, $false, $false, $true, 1, $false, "56-53=3", 2) | Out-Null
$d.Content.Find.Execute("22+65=87", $true, $false, $false, $false, $false, $true, 1, $false, "89-80=9", 2) | Out-Null
$d.Content.Find.Execute("84-6=78", $true, $false, $false, $false, $false, $true, 1, $false, "43+0=43", 2) | Out-Null
$d.Content.Find.Execute("59+28=87", $true, $false, $false, $false, $false, $true, 1, $false, "42-24=18", 2) | Out-Null
$d.Content.Find.Execute("22-14=8", $true, $false, $false, $false, $false, $true, 1, $false, "68+20=88", 2) | Out-Null
$d.Content.Find.Execute("10+66=76", $true, $false, $false, $false, $false, $true, 1, $false, "16+71=87", 2) | Out-Null
$d.Content.Find.Execute("57+17=74", $true, $false, $false, $false, $false, $true, 1, $false, "72+5=77", 2) | Out-Null
$d.Content.Find.Execute("82-26=56", $true, $false, $false, $false, $false, $true, 1, $false, "21+76=97", 2) | Out-Null
$d.Content.Find.Execute("25+24=49", $true, $false, $false, $false, $false, $true, 1, $false, "36+4=40", 2) | Out-Null
$d.Content.Find.Execute("9+69=78", $true, $false, $false, $false, $false, $true, 1, $false, "87-72=15", 2) | Out-Null
$d.Content.Find.Execute("63-17=46", $true, $false, $false, $false, $false, $true, 1, $false, "66-5=61", 2) | Out-Null
$d.Content.Find.Execute("77+13=90", $true, $false, $false, $false, $false, $true, 1, $false, "21-13=8", 2) | Out-Null
$d.Content.Find.Execute("48+47=95", $true, $false, $false, $false, $false, $true, 1, $false, "49+39=88", 2) | Out-Null
$d.Content.Find.Execute("72-6=66", $true, $false, $false, $false, $false, $true, 1, $false, "45+23=68", 2) | Out-Null
$d.Content.Find.Execute("41+9=50", $true, $false, $false, $false, $false, $true, 1, $false, "12-12=0", 2) | Out-Null
$d.Content.Find.Execute("23+60=83", $true, $false, $false, $false, $false, $true, 1, $false, "26+30=56", 2) | Out-Null
$d.Content.Find.Execute("8+61=69", $true, $false, $false, $false, $false, $true, 1, $false, "81+4=85", 2) | Out-Null
$d.Content.Find.Execute("90-87=3", $true, $false, $false, $false, $false, $true, 1, $false, "82-74=8", 2) | Out-Null
$d.Content.Find.Execute("57+15=72", $true, $false, $false, $false, $false, $true, 1, $false, "13+75=88", 2) | Out-Null
$d.Content.Find.Execute("59+15=74", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=55", 2) | Out-Null
$d.Content.Find.Execute("57+1=58", $true, $false, $false, $false, $false, $true, 1, $false, "52+21=73", 2) | Out-Null
$d.Content.Find.Execute("37+52=89", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=57", 2) | Out-Null
$d.Content.Find.Execute("74-4=70", $true, $false, $false, $false, $false, $true, 1, $false, "84-28=56", 2) | Out-Null
$d.Content.Find.Execute("10+3=13", $true, $false, $false, $false, $false, $true, 1, $false, "82-17=65", 2) | Out-Null
$d.Content.Find.Execute("44+24=68", $true, $false, $false, $false, $false, $true, 1, $false, "20+8=28", 2) | Out-Null
$d.Content.Find.Execute("65-0=65", $true, $false, $false, $false, $false, $true, 1, $false, "10+75=85", 2) | Out-Null
$d.Content.Find.Execute("90-43=47", $true, $false, $false, $false, $false, $true, 1, $false, "91-87=4", 2) | Out-Null
$d.Content.Find.Execute("4+50=54", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2) | Out-Null
$d.Content.Find.Execute("86-34=52", $true, $false, $false, $false, $false, $true, 1, $false, "67+23=90", 2) | Out-Null
$d.Content.Find.Execute("55+5=60", $true, $false, $false, $false, $false, $true, 1, $false, "41+28=69", 2) | Out-Null
$d.Content.Find.Execute("11+81=92", $true, $false, $false, $false, $false, $true, 1, $false, "65-13=52", 2) | Out-Null
$d.Content.Find.Execute("36+22=58", $true, $false, $false, $false, $false, $true, 1, $false, "71-31=40", 2) | Out-Null
$d.Content.Find.Execute("18+13=31", $true, $false, $false, $false, $false, $true, 1, $false, "18+30=48", 2) | Out-Null
$d.Content.Find.Execute("39+52=91", $true, $false, $false, $false, $false, $true, 1, $false, "79-28=51", 2) | Out-Null
$d.Content.Find.Execute("39+15=54", $true, $false, $false, $false, $false, $true, 1, $false, "90-86=4", 2) | Out-Null
$d.Content.Find.Execute("28+45=73", $true, $false, $false, $false, $false, $true, 1, $false, "11-6=5", 2) | Out-Null
$d.Content.Find.Execute("35+37=72", $true, $false, $false, $false, $false, $true, 1, $false, "37+12=49", 2) | Out-Null
$d.Content.Find.Execute("13+19=32", $true, $false, $false, $false, $false, $true, 1, $false, "23-3=20", 2) | Out-Null
$d.Content.Find.Execute("98-95=3", $true, $false, $false, $false, $false, $true, 1, $false, "13+70=83", 2) | Out-Null
$d.Content.Find.Execute("28-15=13", $true, $false, $false, $false, $false, $true, 1, $false, "88-31=57", 2) | Out-Null
$d.Content.Find.Execute("99-94=5", $true, $false, $false, $false, $false, $true, 1, $false, "23-13=10", 2) | Out-Null
$d.Content.Find.Execute("54-33=21", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=17", 2) | Out-Null
$d.Content.Find.Execute("39+21=60", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=15", 2) | Out-Null
$d.Content.Find.Execute("52-9=43", $true, $false, $false, $false, $false, $true, 1, $false, "42+22=64", 2) | Out-Null
$d.Content.Find.Execute("66-8=58", $true, $false, $false, $false, $false, $true, 1, $false, "12+59=71", 2) | Out-Null
$d.Content.Find.Execute("90-51=39", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("49-15=34", $true, $false, $false, $false, $false, $true, 1, $false, "5+3=8", 2) | Out-Null
$d.Content.Find.Execute("53-16=37", $true, $false, $false, $false, $false, $true, 1, $false, "57+26=83", 2) | Out-Null
$d.Content.Find.Execute("18+60=78", $true, $false, $false, $false, $false, $true, 1, $false, "32-13=19", 2) | Out-Null
$d.Content.Find.Execute("70-49=21", $true, $false, $false, $false, $false, $true, 1, $false, "0+73=73", 2) | Out-Null
$d.Content.Find.Execute("2+49=51", $true, $false, $false, $false, $false, $true, 1, $false, "37-7=30", 2) | Out-Null
$d.Content.Find.Execute("65+34=99", $true, $false, $false, $false, $false, $true, 1, $false, "43+39=82", 2) | Out-Null
$d.Content.Find.Execute("87-55=32", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=12", 2) | Out-Null
$d.Content.Find.Execute("50-6=44", $true, $false, $false, $false, $false, $true, 1, $false, "85-11=74", 2) | Out-Null
